$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new "StatQuery" text (column C) used to replace the old stat query for
# every row (CasesTab / SamplesTab / FilesTab). Built as a here-string so the
# embedded tabs, single/double quotes and blank lines survive verbatim.
$newStatQuery = @"
MATCH (p:program)<--(s:study)<--(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
      WHERE (size([]) = 0 OR s.clinical_study_designation IN [])
        AND (s.study_disposition = 'Unrestricted')
        AND s.clinical_study_designation IN ['UBC01']
		 and diag.stage_of_disease in [ 'T2N0M0', 'T2N1M0', 'T3N0M0', 'T3N0M1', 'T3N1M0', 'Not Applicable'] 
        AND (size([]) = 0 OR demo.sex IN [])
        AND (size([]) = 0 OR demo.neutered_indicator IN [])
        AND (size([]) = 0 OR diag.disease_term IN [])
        AND (size([]) = 0 OR diag.primary_disease_site IN [])
        AND (size([]) = 0 OR diag.stage_of_disease IN [])
        AND (size([]) = 0 OR diag.best_response IN [])
    OPTIONAL MATCH (c)-->(co:cohort)
    OPTIONAL MATCH (f:file)-[*]->(c)
    OPTIONAL MATCH (f)-->(parent)
    OPTIONAL MATCH (samp:sample)-->(c)
    OPTIONAL MATCH (samp)<--(al:aliquot)
    WITH DISTINCT c AS c, p, s, co, demo, diag, f, parent, samp, al
      WHERE (size([]) = 0 OR samp.summarized_sample_type IN [])
        AND (size([]) = 0 OR samp.specific_sample_pathology IN [])
        AND (size([]) = 0 OR samp.sample_site IN [])
        AND (size([]) = 0 OR head(labels(parent)) IN [])
        AND (size([]) = 0 OR f.file_type IN [])
        AND (size([]) = 0 OR f.file_format IN [])
    WITH c.case_id AS case_id,
         s.clinical_study_designation AS study_code,
         s.clinical_study_type AS study_type,
         co.cohort_description AS cohort,
         demo.breed AS breed,
         diag.disease_term AS diagnosis,
         diag.stage_of_disease AS stage_of_disease,
         diag.primary_disease_site AS disease_site,
         demo.patient_age_at_enrollment AS age,
         demo.sex AS sex,
         demo.neutered_indicator AS neutered_status,
         demo.weight AS weight,
         diag.best_response AS response_to_treatment,
         samp.sample_id AS sample_id,
         f.uuid AS file_id,
         al
    RETURN
COUNT(DISTINCT file_id) as number_of_files,
COUNT(DISTINCT sample_id) as number_of_sample,
COUNT(DISTINCT case_id) as number_of_cases,
COUNT(DISTINCT study_code) as number_of_study,
COUNT(DISTINCT al) as number_of_aliquot
    
"@

# Update the StatQuery column (C) for the three data rows: CasesTab (2),
# SamplesTab (3), FilesTab (4). The "query" column (B) text is unchanged.
$ws.Range("C2").Value2 = $newStatQuery
$ws.Range("C3").Value2 = $newStatQuery
$ws.Range("C4").Value2 = $newStatQuery

# Update the view/selection to match the post-edit state (scrolled to row 4,
# active cell C4).
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("C4").Select()
